# Apply IGCC Netting Flows data update:
#  - Shift data window from day 45984 to day 45993 (9 days later)
#  - Extend data from 32 rows (rows 2-33) to 42 rows (rows 2-43)
#  - Update Import (col B) / Export (col C) values accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45993, 1.39, 14.36),
    @(45993.01041666666, 0, 1.055),
    @(45993.02083333334, 0, 3.328),
    @(45993.03125, 0, 4.428),
    @(45993.04166666666, 0, 1.424),
    @(45993.05208333334, 0, 18.278),
    @(45993.0625, 0, 3.819),
    @(45993.07291666666, 0, 5.415),
    @(45993.08333333334, 0, 0.569),
    @(45993.09375, 0, 0.639),
    @(45993.10416666666, 4.973, 0.013),
    @(45993.11458333334, 0.002, 0.785),
    @(45993.125, 0, 8.766999999999999),
    @(45993.13541666666, 0.042, 0.536),
    @(45993.14583333334, 0.159, 0.713),
    @(45993.15625, 2.944, 0.416),
    @(45993.16666666666, 4.127, 0.191),
    @(45993.17708333334, 5.369, 0.537),
    @(45993.1875, 1.54, 1.903),
    @(45993.19791666666, 3.79, 0.617),
    @(45993.20833333334, 12.263, 2.751),
    @(45993.21875, 10.276, 0.173),
    @(45993.22916666666, 10.383, 0.081),
    @(45993.23958333334, 8.279, 0.061),
    @(45993.25, 7.695, 0),
    @(45993.26041666666, 14.984, 0),
    @(45993.27083333334, 28.921, 0),
    @(45993.28125, 19.787, 0.004),
    @(45993.29166666666, 1.629, 1.951),
    @(45993.30208333334, 1.933, 0),
    @(45993.3125, 2.433, 0),
    @(45993.32291666666, 0.846, 7.175),
    @(45993.33333333334, 0, 39.747),
    @(45993.34375, 0, 23.315),
    @(45993.35416666666, 1.455, 1.32),
    @(45993.36458333334, 3.101, 7.265),
    @(45993.375, 0.134, 1.604),
    @(45993.38541666666, 0.157, 4.605),
    @(45993.39583333334, 4.15, 0),
    @(45993.40625, 0.963, 0.475),
    @(45993.41666666666, 9.005000000000001, 0),
    @(45993.42708333334, 0, 0)
)

$numRows = $data.Length
for ($i = 0; $i -lt $numRows; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}

# Newly added rows (34-43) need the same timestamp number format as the
# pre-existing rows in column A (YYYY-MM-DD HH:MM:SS).
$ws.Range("A2:A43").NumberFormat = "YYYY-MM-DD HH:MM:SS"
